# The deck's two theme parts ("Office Theme" colours and the "Integral"
# colours that are actually wired up as the presentation's live theme)
# had their colour values swapped: the live theme (reached through any
# slide's ThemeColorScheme, which is backed by the shared theme part
# used by the slide master / the whole deck) picks up the plain
# "Office Theme" palette, while the palette that used to live there
# ("Integral") is what the notes master's theme now carries.
#
# Modifying ThemeColorScheme.Item(n).RGB on a slide edits the colour
# values inside the shared <a:clrScheme> of the presentation's real
# theme part in place (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -
# in that index order), which is exactly the content that changed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      (BGR order = 0x44546A)
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      (BGR order = 0xE7E6E6)
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  (BGR order = 0x5B9BD5)
$tcs.Item(6).RGB  = 0x317DED   # accent2  (BGR order = 0xED7D31)
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  (BGR order = 0xA5A5A5)
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  (BGR order = 0xFFC000)
$tcs.Item(9).RGB  = 0xC47244   # accent5  (BGR order = 0x4472C4)
$tcs.Item(10).RGB = 0x47AD70   # accent6  (BGR order = 0x70AD47)
$tcs.Item(11).RGB = 0xC16305   # hlink    (BGR order = 0x0563C1)
$tcs.Item(12).RGB = 0x724F95   # folHlink (BGR order = 0x954F72)
